$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column from 2023-11-13 to 2023-11-14
# for rows 2 through 7 (serial date 45243 -> 45244).
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45244
}
